$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "to do" notes under a new section (rows 30-33), then add the
# new section header "Katelyn Gone Now" at A29 last so it claims the final
# (highest-index) shared-string slot, matching authoring order where the
# header was styled/typed after the bullet notes beneath it.
$ws.Range("A30").Value = "Dewpoint gives error with 0 humidity (likely a divide by zero issue)"
$ws.Range("A31").Value = "Current data likely needs to be two separate dictionaries: when an error occurs, one sensor is pulling data from currentdata and that could actually be the data of the other sensor"
$ws.Range("A32").Value = "Check in arduino that it isn't sending 255 from a sensor: make anything over 255 into 254"
$ws.Range("A33").Value = "X axis… put time in here"

$ws.Range("A29").Value = "Katelyn Gone Now"
$ws.Range("A29").Style = "Heading 2"

# Match the author's final selection position.
$ws.Range("F31").Select() | Out-Null
